$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 0.799958874454461
$ws.Range("C2").Value2 = 0.1067715799096476
$ws.Range("E2").Value2 = 0.1046706745146118
$ws.Range("F2").Value2 = 0.4443680307746263
$ws.Range("G2").Value2 = 1.63393404254073
$ws.Range("H2").Value2 = 1.455794365372839
$ws.Range("K2").Value2 = 0.4078566521046696
$ws.Range("L2").Value2 = 0.2041937060047871
$ws.Range("B3").Value2 = 0.7641185369060679
$ws.Range("C3").Value2 = 0.1059563866660902
$ws.Range("E3").Value2 = 0.1038750138918836
$ws.Range("F3").Value2 = 0.387822817061874
$ws.Range("G3").Value2 = 1.629867427624916
$ws.Range("H3").Value2 = 1.459036684422259
$ws.Range("K3").Value2 = 0.376563955205711
$ws.Range("L3").Value2 = 0.1976209913672307
$ws.Range("B4").Value2 = 0.7425441335375353
$ws.Range("C4").Value2 = 0.1054441480863524
$ws.Range("E4").Value2 = 0.1034352773344018
$ws.Range("F4").Value2 = 0.3531389305168915
$ws.Range("G4").Value2 = 1.628179123928149
$ws.Range("H4").Value2 = 1.461583905893761
$ws.Range("K4").Value2 = 0.3575402849988194
$ws.Range("L4").Value2 = 0.1937017606587119
$ws.Range("B5").Value2 = 0.733861256288435
$ws.Range("C5").Value2 = 0.10523246273452
$ws.Range("E5").Value2 = 0.1032683607074141
$ws.Range("F5").Value2 = 0.3390132514313251
$ws.Range("G5").Value2 = 1.627694002172944
$ws.Range("H5").Value2 = 1.462761735369867
$ws.Range("K5").Value2 = 0.349836000203041
$ws.Range("L5").Value2 = 0.1921339338030918
$ws.Range("B6").Value2 = 0.7324260557452078
$ws.Range("C6").Value2 = 0.1051971347308758
$ws.Range("E6").Value2 = 0.1032413862217254
$ws.Range("F6").Value2 = 0.336668177824194
$ws.Range("G6").Value2 = 1.6276256878959
$ws.Range("H6").Value2 = 1.462965754647584
$ws.Range("K6").Value2 = 0.3485596154509381
$ws.Range("L6").Value2 = 0.1918753668341964
$ws.Range("B7").Value2 = 0.7424265920029995
$ws.Range("C7").Value2 = 0.1054413051432839
$ws.Range("E7").Value2 = 0.1034329765080102
$ws.Range("F7").Value2 = 0.3529483938344953
$ws.Range("G7").Value2 = 1.628171760582532
$ws.Range("H7").Value2 = 1.461599224545452
$ws.Range("K7").Value2 = 0.3574361876456891
$ws.Range("L7").Value2 = 0.1936804977811732
$ws.Range("B8").Value2 = 0.787511696376157
$ws.Range("C8").Value2 = 0.1064929306021547
$ws.Range("E8").Value2 = 0.1043862060051275
$ws.Range("F8").Value2 = 0.4248636149813336
$ws.Range("G8").Value2 = 1.632363766450908
$ws.Range("H8").Value2 = 1.456796768202892
$ws.Range("K8").Value2 = 0.3970275649936639
$ws.Range("L8").Value2 = 0.2019032701317656
$ws.Range("B9").Value2 = 0.8793409728630763
$ws.Range("C9").Value2 = 0.1084623716825064
$ws.Range("E9").Value2 = 0.1066425885496294
$ws.Range("F9").Value2 = 0.5661985755041457
$ws.Range("G9").Value2 = 1.64702298814268
$ws.Range("H9").Value2 = 1.451799630325709
$ws.Range("K9").Value2 = 0.4761710763863221
$ws.Range("L9").Value2 = 0.2189526822665044
$ws.Range("B10").Value2 = 0.9488892631993622
$ws.Range("C10").Value2 = 0.1098529993078117
$ws.Range("E10").Value2 = 0.1085365086331542
$ws.Range("F10").Value2 = 0.6702781546542269
$ws.Range("G10").Value2 = 1.661752823435734
$ws.Range("H10").Value2 = 1.450831983076966
$ws.Range("K10").Value2 = 0.5352366695822184
$ws.Range("L10").Value2 = 0.2320451193691468
$ws.Range("B11").Value2 = 0.9809807077458288
$ws.Range("C11").Value2 = 0.1104734654128166
$ws.Range("E11").Value2 = 0.1094494405650224
$ws.Range("F11").Value2 = 0.7176906081379002
$ws.Range("G11").Value2 = 1.669321225928627
$ws.Range("H11").Value2 = 1.45098089062526
$ws.Range("K11").Value2 = 0.5623074566917126
$ws.Range("L11").Value2 = 0.2381248303260435
$ws.Range("B12").Value2 = 0.9931979675714615
$ws.Range("C12").Value2 = 0.1107066771942726
$ws.Range("E12").Value2 = 0.1098025306382127
$ws.Range("F12").Value2 = 0.7356546913071611
$ws.Range("G12").Value2 = 1.672312515796932
$ws.Range("H12").Value2 = 1.451122133191461
$ws.Range("K12").Value2 = 0.5725873678381674
$ws.Range("L12").Value2 = 0.2404448949145603
$ws.Range("B13").Value2 = 0.9905638791410638
$ws.Range("C13").Value2 = 0.110656528482572
$ws.Range("E13").Value2 = 0.1097261581192335
$ws.Range("F13").Value2 = 0.7317853510981394
$ws.Range("G13").Value2 = 1.671662706937155
$ws.Range("H13").Value2 = 1.451087937860137
$ws.Range("K13").Value2 = 0.570372127525701
$ws.Range("L13").Value2 = 0.2399444354789324
$ws.Range("B14").Value2 = 0.9819845295315872
$ws.Range("C14").Value2 = 0.110492686863239
$ws.Range("E14").Value2 = 0.1094783415716698
$ws.Range("F14").Value2 = 0.7191683204515869
$ws.Range("G14").Value2 = 1.669564806957197
$ws.Range("H14").Value2 = 1.450990809302965
$ws.Range("K14").Value2 = 0.5631526152227764
$ws.Range("L14").Value2 = 0.2383153465293759
$ws.Range("B15").Value2 = 0.9767378767456307
$ws.Range("C15").Value2 = 0.110392101860306
$ws.Range("E15").Value2 = 0.1093275080455349
$ws.Range("F15").Value2 = 0.7114413442032514
$ws.Range("G15").Value2 = 1.668296115901626
$ws.Range("H15").Value2 = 1.450942370080952
$ws.Range("K15").Value2 = 0.5587341993902157
$ws.Range("L15").Value2 = 0.2373198017288303
$ws.Range("B16").Value2 = 0.9468011133989194
$ws.Range("C16").Value2 = 0.1098122061557376
$ws.Range("E16").Value2 = 0.1084778797674524
$ws.Range("F16").Value2 = 0.6671810134426437
$ws.Range("G16").Value2 = 1.661275716050341
$ws.Range("H16").Value2 = 1.450834117959801
$ws.Range("K16").Value2 = 0.5334715692610246
$ws.Range("L16").Value2 = 0.2316502885245768
$ws.Range("B17").Value2 = 0.9285518373754087
$ws.Range("C17").Value2 = 0.109453350241381
$ws.Range("E17").Value2 = 0.1079698155364319
$ws.Range("F17").Value2 = 0.6400460337125793
$ws.Range("G17").Value2 = 1.657191545362025
$ws.Range("H17").Value2 = 1.450918692029518
$ws.Range("K17").Value2 = 0.5180252486725578
$ws.Range("L17").Value2 = 0.2282039539361307
$ws.Range("B18").Value2 = 0.9180980468719326
$ws.Range("C18").Value2 = 0.1092458033690917
$ws.Range("E18").Value2 = 0.1076824267132039
$ws.Range("F18").Value2 = 0.6244449056556647
$ws.Range("G18").Value2 = 1.654924081740759
$ws.Range("H18").Value2 = 1.451022776596545
$ws.Range("K18").Value2 = 0.5091599158341751
$ws.Range("L18").Value2 = 0.2262333754005112
$ws.Range("B19").Value2 = 0.9145659177460459
$ws.Range("C19").Value2 = 0.1091753353482687
$ws.Range("E19").Value2 = 0.1075859526031024
$ws.Range("F19").Value2 = 0.619163680173358
$ws.Range("G19").Value2 = 1.654170363931783
$ws.Range("H19").Value2 = 1.451067535330594
$ws.Range("K19").Value2 = 0.5061615325871003
$ws.Range("L19").Value2 = 0.2255681745165958
$ws.Range("B20").Value2 = 0.930490085772874
$ws.Range("C20").Value2 = 0.109491669305811
$ws.Range("E20").Value2 = 0.1080233993447024
$ws.Range("F20").Value2 = 0.642933953830422
$ws.Range("G20").Value2 = 1.657617858736216
$ws.Range("H20").Value2 = 1.450903950120477
$ws.Range("K20").Value2 = 0.5196675726817261
$ws.Range("L20").Value2 = 0.2285696152134165
$ws.Range("B21").Value2 = 0.9845027327428113
$ws.Range("C21").Value2 = 0.1105408584831196
$ws.Range("E21").Value2 = 0.1095509309820173
$ws.Range("F21").Value2 = 0.7228739723491628
$ws.Range("G21").Value2 = 1.670177606267657
$ws.Range("H21").Value2 = 1.451017034242398
$ws.Range("K21").Value2 = 0.5652723810962641
$ws.Range("L21").Value2 = 0.2387933661301673
$ws.Range("B22").Value2 = 1.02018142565521
$ws.Range("C22").Value2 = 0.1112163916819995
$ws.Range("E22").Value2 = 0.1105922899407759
$ws.Range("F22").Value2 = 0.7751780083420101
$ws.Range("G22").Value2 = 1.679116714184858
$ws.Range("H22").Value2 = 1.451585605181805
$ws.Range("K22").Value2 = 0.5952454935945468
$ws.Range("L22").Value2 = 0.2455789863393534
$ws.Range("B23").Value2 = 1.001104517677049
$ws.Range("C23").Value2 = 0.110856777617542
$ws.Range("E23").Value2 = 0.110032561741658
$ws.Range("F23").Value2 = 0.7472568307830727
$ws.Range("G23").Value2 = 1.674278722934105
$ws.Range("H23").Value2 = 1.451236838895881
$ws.Range("K23").Value2 = 0.5792329937249008
$ws.Range("L23").Value2 = 0.2419478765084904
$ws.Range("B24").Value2 = 0.9296136854048029
$ws.Range("C24").Value2 = 0.1094743491069892
$ws.Range("E24").Value2 = 0.107999159452369
$ws.Range("F24").Value2 = 0.6416283278902171
$ws.Range("G24").Value2 = 1.657424871487649
$ws.Range("H24").Value2 = 1.450910442177275
$ws.Range("K24").Value2 = 0.5189250313773357
$ws.Range("L24").Value2 = 0.2284042662194565
$ws.Range("B25").Value2 = 0.854133089492251
$ws.Range("C25").Value2 = 0.1079394948699779
$ws.Range("E25").Value2 = 0.1059907183835378
$ws.Range("F25").Value2 = 0.5279251897347166
$ws.Range("G25").Value2 = 1.642364404426402
$ws.Range("H25").Value2 = 1.452677327656133
$ws.Range("K25").Value2 = 0.4545994184190079
$ws.Range("L25").Value2 = 0.2142411172982861
